$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Protocol text corrections ---
# "Ticket de lotterie gagnant" had a typo ("lotterie" -> "loterie"); fix every
# occurrence (D7 and C13 in the original layout).
$ws.Range("D7").Value = "Ticket de loterie gagnant"
$ws.Range("C13").Value = "Ticket de loterie gagnant"

# C7 used to share the "Clé de voiture" string with B13; give tour 1 its own,
# more distinctive prize so it isn't a duplicate of the "tour 7" reward.
$ws.Range("C7").Value = "Clé de voiture de luxe"

# --- Restore the cursor position that was left selected when the file was saved ---
$ws.Range("D15").Select() | Out-Null
